$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'72.304.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.62%  "
$ws.Range("D3").Value = "'4.054.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.81%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'526.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").Value = "'148.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.31%  "
$ws.Range("D7").Value = "'0.714"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +18.95%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.765"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.15%  "
$ws.Range("E10").Value = "  +7.62%  "
$ws.Range("D11").Value = "'0.0000334"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.96%  "
$ws.Range("D12").Value = "'48.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +19.33%  "
$ws.Range("D13").Value = "'10.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.53%  "
$ws.Range("D14").Value = "'4.702.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.84%  "
$ws.Range("D15").Value = "'4.037.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.90%  "
$ws.Range("D16").Value = "'14.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Value = "'21.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.38%  "
$ws.Range("E18").Value = "  +3.14%  "
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").Value = "'72.230.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.68%  "
$ws.Range("D21").Value = "'439.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.37%  "
$ws.Range("D22").Value = "'101.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +18.51%  "
$ws.Range("D23").Value = "'3.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.07%  "
$ws.Range("D24").Value = "'14.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.51%  "
$ws.Range("D25").Value = "'4.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.30%  "
$ws.Range("D26").Value = "'11.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.88%  "
$ws.Range("D27").Value = "'11.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.01%  "
$ws.Range("D28").Value = "'37.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.53%  "
$ws.Range("D29").Value = "'5.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.10%  "
$ws.Range("D30").Value = "'3.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +21.70%  "
$ws.Range("D31").Value = "'13.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.19%  "
$ws.Range("E32").Value = "  +8.03%  "
$ws.Range("D33").Value = "'679.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("E34").Value = "  +10.35%  "
$ws.Range("D35").Value = "'66.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.25%  "
$ws.Range("D36").Value = "'42.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.44%  "
$ws.Range("D39").Value = "'0.157"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.43%  "
$ws.Range("D40").Value = "'3.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.99%  "
$ws.Range("D41").Value = "'0.0510"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.12%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "'3.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").Value = "'0.156"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +13.61%  "
$ws.Range("D46").Value = "'2.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("D47").Value = "'3.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("D48").Value = "'9.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +13.04%  "
$ws.Range("E49").Value = "  +7.44%  "

# Row 37: PEPE -> TheGraph
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").Value = "'0.438"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.96%  "

# Row 38: TheGraph -> PEPE
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0870"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.55%  "

# Row 50: FLOKI -> LidoDAOToken
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").Value = "'3.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.80%  "

# Row 51: LidoDAOToken -> FLOKI
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "'0.000275"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.03%  "
